$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, applied as text (preserving formatting
# such as trailing zeros and thousand-dot separators, and avoiding Excel
# auto-converting numeric-looking strings into floating point numbers).
$updates = [ordered]@{
    "D2" = "64.630.29"
    "E2" = "  -0.21%  "
    "D3" = "3.341.30"
    "E3" = "  -0.91%  "
    "E4" = "  +0.00%  "
    "D5" = "552.84"
    "E5" = "  -0.48%  "
    "D6" = "173.53"
    "E6" = "  -1.60%  "
    "D7" = "0.627"
    "E7" = "  +1.71%  "
    "B8" = "LidoStakedEther"
    "C8" = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
    "D8" = "3.332.96"
    "E8" = "  -0.92%  "
    "B9" = "USDC"
    "C9" = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
    "D9" = "1.00"
    "E9" = "  +0.02%  "
    "D10" = "0.174"
    "E10" = "  +6.85%  "
    "D11" = "0.638"
    "E11" = "  +1.53%  "
    "D12" = "53.44"
    "E12" = "  -3.26%  "
    "D13" = "0.0000279"
    "E13" = "  +2.20%  "
    "D14" = "9.09"
    "E14" = "  +0.20%  "
    "D15" = "3.863.90"
    "E15" = "  -1.01%  "
    "E16" = "  +1.90%  "
    "D17" = "18.14"
    "E17" = "  -1.30%  "
    "D18" = "3.347.91"
    "E18" = "  -0.41%  "
    "D19" = "64.437.45"
    "E19" = "  -0.19%  "
    "D20" = "11.74"
    "E20" = "  -0.82%  "
    "D21" = "0.987"
    "E21" = "  +0.56%  "
    "D22" = "447.42"
    "E22" = "  +2.41%  "
    "D23" = "4.96"
    "E23" = "  -0.67%  "
    "D24" = "4.05"
    "E24" = "  -0.92%  "
    "D25" = "87.11"
    "E25" = "  +3.11%  "
    "D26" = "13.74"
    "E26" = "  +3.48%  "
    "D27" = "2.87"
    "E27" = "  +0.88%  "
    "D28" = "10.64"
    "E28" = "  -1.86%  "
    "D29" = "8.60"
    "E29" = "  -2.15%  "
    "D30" = "30.89"
    "E30" = "  +3.61%  "
    "D31" = "6.52"
    "E31" = "  -1.87%  "
    "B32" = "OKB"
    "C32" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D32" = "62.54"
    "E32" = "  +6.60%  "
    "B33" = "Cosmos"
    "C33" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D33" = "11.39"
    "E33" = "  -0.82%  "
    "D34" = "571.36"
    "E34" = "  -1.84%  "
    "D35" = "0.107"
    "E35" = "  -1.32%  "
    "E36" = "  +0.11%  "
    "E37" = "  -1.24%  "
    "D38" = "3.55"
    "E38" = "  +0.15%  "
    "D39" = "35.35"
    "E39" = "  -1.31%  "
    "D40" = "0.368"
    "E40" = "  -0.30%  "
    "D41" = "0.0₃0733"
    "E41" = "  -3.54%  "
    "D42" = "3.071.82"
    "E42" = "  -1.46%  "
    "D43" = "0.0415"
    "E43" = "  +1.07%  "
    "E44" = "  -3.27%  "
    "B45" = "ApeXProtocol"
    "C45" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "D45" = "3.21"
    "E45" = "  -2.61%  "
    "B46" = "Stellar"
    "C46" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D46" = "0.134"
    "E46" = "  +3.14%  "
    "E47" = "  -1.65%  "
    "D48" = "0.998"
    "E48" = "  +0.07%  "
    "D49" = "140.61"
    "E49" = "  +4.40%  "
    "E50" = "  -3.00%  "
    "D51" = "8.23"
    "E51" = "  -1.17%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force a Text number format before assigning so Excel stores the
    # value as a string instead of silently coercing it to a Double
    # (which would corrupt values like "1.00" -> 1 or "552.84" -> 552.84000000000003).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Drop the temporary Text format again so the cell keeps the workbook
    # default style (no explicit "s" attribute), matching the original file.
    $cell.ClearFormats()
}
